$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.538.66'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.23%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.876.70'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.59%  '

$ws.Range("E4").Value = '  +0.65%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.54'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.23%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.014'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.63%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4795'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.49%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3778'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.48%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07391'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.69%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9409'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.55%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.78'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +5.85%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07871'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.35%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.882.66'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.71%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.450'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.71%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.614'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.29%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '91.18'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.99%  '

$ws.Range("E17").Value = '  +0.61%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008998'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +4.05%  '

$ws.Range("E19").Value = '  +0.58%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.98'
$ws.Range("D20").Style = "Normal"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '27.563.28'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.20%  '

$ws.Range("E22").Value = '  +2.42%  '

$ws.Range("E23").Value = '  +0.62%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.954'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.87%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '153.98'
$ws.Range("D25").Style = "Normal"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '18.58'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.23%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.027'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.89%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '116.18'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.61%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.004'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.82%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.08927'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.63%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.322'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.37%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.219'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.58%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.613'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.64%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7527'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.58%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.699'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.15%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02078'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +6.64%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.122'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.85%  '

$ws.Range("E38").Value = '  +0.94%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.003'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.01%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5372'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.85%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.097'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.79%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1521'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.51%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.451'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.90%  '

$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.64'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.96%  '

$ws.Range("B45").Value = 'Decentraland'
$ws.Range("C45").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4847'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.05%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.015'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.65%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.665'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.70%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '103.16'
$ws.Range("D48").Style = "Normal"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '67.19'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.48%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06112'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.23%  '

$ws.Range("E51").Value = '  +1.86%  '
